$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts D:K -> E:L), preserving values/styles.
$ws.Columns("D:D").Insert()

# Copy number formats from column E (the old column D, now shifted) onto the
# newly inserted column D so the new cells pick up the correct style (date
# format on the header rows, integer format on the data rows).
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the new quarter's figures (period ending
# 2018-09-30, serial 43373).
$ws.Range("D7").Value = 43373
$ws.Range("D8").Value = 21000
$ws.Range("D9").Value = 18800
$ws.Range("D10").Value = 2200
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 20400
$ws.Range("D18").Value = 600
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = 600
$ws.Range("D22").Value = 100
$ws.Range("D23").Value = 500
$ws.Range("D24").Value = 200
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 300
$ws.Range("D27").Value = 200
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = 200
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 200
$ws.Range("D38").Value = 43373
$ws.Range("D41").Value = 1000
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 11300
$ws.Range("D44").Value = 11400
$ws.Range("D45").Value = 4700
$ws.Range("D46").Value = 28300
$ws.Range("D47").Value = 1800
$ws.Range("D48").Value = 200
$ws.Range("D49").Value = 600
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 100
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 31000
$ws.Range("D57").Value = 15600
$ws.Range("D58").Value = 5700
$ws.Range("D59").Value = 3500
$ws.Range("D60").Value = 24800
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 25800
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 200
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 5200
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43373
$ws.Range("D81").Value = 200
$ws.Range("D83").Value = 100
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 1300
$ws.Range("D91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 0
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -1700
$ws.Range("D101").Value = -100
$ws.Range("D102").Value = -500
